$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking values
# (e.g. "581.95", "62.139.67") are stored as strings, not numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '62.139.67'
$ws.Range("E2").Value = '  -0.62%  '

$ws.Range("D3").Value = '2.445.99'
$ws.Range("E3").Value = '  +0.26%  '

$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").Value = '581.95'
$ws.Range("E5").Value = '  +1.36%  '

$ws.Range("D6").Value = '142.78'
$ws.Range("E6").Value = '  -1.01%  '

$ws.Range("E7").Value = '  +0.10%  '

$ws.Range("D8").Value = '0.533'
$ws.Range("E8").Value = '  +0.33%  '

$ws.Range("D9").Value = '2.440.86'
$ws.Range("E9").Value = '  +0.18%  '

$ws.Range("D10").Value = '0.110'
$ws.Range("E10").Value = '  +1.82%  '

$ws.Range("E11").Value = '  +2.37%  '

$ws.Range("E13").Value = '  -2.56%  '

$ws.Range("D14").Value = '26.39'
$ws.Range("E14").Value = '  -0.73%  '

$ws.Range("E15").Value = '  -0.03%  '

$ws.Range("D16").Value = '2.880.58'
$ws.Range("E16").Value = '  +0.16%  '

$ws.Range("D17").Value = '62.060.25'
$ws.Range("E17").Value = '  -0.40%  '

$ws.Range("D18").Value = '2.438.12'
$ws.Range("E18").Value = '  +0.90%  '

$ws.Range("E19").Value = '  -3.78%  '

$ws.Range("D20").Value = '7.21'
$ws.Range("E20").Value = '  +0.65%  '

$ws.Range("D21").Value = '326.07'
$ws.Range("E21").Value = '  -0.55%  '

$ws.Range("D22").Value = '4.10'
$ws.Range("E22").Value = '  -1.27%  '

$ws.Range("E23").Value = '  +0.06%  '

$ws.Range("E24").Value = '  -6.07%  '

$ws.Range("D25").Value = '65.59'
$ws.Range("E25").Value = '  +0.13%  '

$ws.Range("D26").Value = '9.12'
$ws.Range("E26").Value = '  -0.68%  '

$ws.Range("D27").Value = '597.36'

$ws.Range("D28").Value = '0.0₃0966'
$ws.Range("E28").Value = '  +0.22%  '

$ws.Range("E30").Value = '  +0.21%  '

$ws.Range("E31").Value = '  -1.49%  '

$ws.Range("D32").Value = '7.98'
$ws.Range("E32").Value = '  -1.21%  '

$ws.Range("D33").Value = '1.90'
$ws.Range("E33").Value = '  +0.85%  '

$ws.Range("E34").Value = '  -0.56%  '

$ws.Range("E35").Value = '  -2.48%  '

$ws.Range("E36").Value = '  +0.21%  '

$ws.Range("E37").Value = '  -1.94%  '

$ws.Range("D38").Value = '0.376'
$ws.Range("E38").Value = '  +0.13%  '

$ws.Range("D39").Value = '154.00'
$ws.Range("E39").Value = '  +5.35%  '

$ws.Range("D40").Value = '18.41'
$ws.Range("E40").Value = '  -0.56%  '

$ws.Range("D41").Value = '5.28'
$ws.Range("E41").Value = '  +0.74%  '

$ws.Range("D42").Value = '43.24'
$ws.Range("E42").Value = '  +2.31%  '

$ws.Range("D43").Value = '1.71'
$ws.Range("E43").Value = '  -1.53%  '

$ws.Range("E44").Value = '  +0.04%  '

$ws.Range("D45").Value = '2.53'
$ws.Range("E45").Value = '  +1.60%  '

$ws.Range("D46").Value = '0.0₆0276'
$ws.Range("E46").Value = '  +23.06%  '

$ws.Range("D47").Value = '141.75'
$ws.Range("E47").Value = '  -2.52%  '

$ws.Range("D48").Value = '3.62'
$ws.Range("E48").Value = '  -2.46%  '

$ws.Range("D49").Value = '0.601'
$ws.Range("E49").Value = '  +0.50%  '

$ws.Range("D50").Value = '0.0518'
$ws.Range("E50").Value = '  -1.14%  '

$ws.Range("D51").Value = '19.82'
$ws.Range("E51").Value = '  +0.68%  '

# Restore default cell style on column D so no residual number-format
# styling is left on the cells (matches original unstyled cells).
$ws.Range("D2:D51").Style = "Normal"
